# Daily attendance processing - 2025-12-28 17:56:39
#
# The "Recorded By" column (G) lists the users who recorded/edited a
# session. For a specific set of rows the literal "System" entry was
# sitting first in the comma-separated list; it should instead be moved
# to the end of the list (the other recorder(s) keep their relative
# order).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targetRows = @(2,4,5,8,11,17,28,30,31,34,37,43,54,56,57,60,63,69,80,81,82,93,94,96,106,107,108,119,120,122,132,133,134,145,146,148)

foreach ($row in $targetRows) {
    $cell = $ws.Cells.Item($row, 7)   # column G
    $current = $cell.Value2

    $parts = $current -split ", "
    $newValue = (($parts[1..($parts.Length - 1)]) + $parts[0]) -join ", "

    $cell.Value = $newValue
}
